# Apply the VWSabbreviations updates:
#  - update a couple of existing definitions / abbreviations
#  - replace the "bottom technique" rows (19-24) with new "top/neutral technique" rows
#  - remove the now-obsolete rows that used to follow them (old rows 25-32)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Focus abbreviation -> "none"
$ws.Range("C4").Value = "none"

# Weight definition gains units
$ws.Range("D6").Value = "Weight for compeition in kilograms"

# Counter abbreviation -> "GB"
$ws.Range("C17").Value = "GB"

# Replace rows 19-24 content (previously Stand up/Escape/Reversal/Cut/Breakdown/Mat Return)
$ws.Range("B19").Value = "Exposure"
$ws.Range("C19").Value = "Expo"
$ws.Range("D19").Value = "Neutral technique for forcing another wrestler to expose their back"

$ws.Range("B20").Value = "Recovery"
$ws.Range("C20").Value = "Recovery"
$ws.Range("D20").Value = "Recovering action after being exposed"

$ws.Range("B21").Value = "Gut"
$ws.Range("C21").Value = "Gut"
$ws.Range("D21").Value = "Top technique focuses on torso"

$ws.Range("B22").Value = "LegLace"
$ws.Range("C22").Value = "LegLace"
$ws.Range("D22").Value = "Top technique focuses on legs"

$ws.Range("B23").Value = "Turn"
$ws.Range("C23").Value = "Turn"
$ws.Range("D23").Value = "Alternative top techniques"

$ws.Range("B24").Value = "Passive"
$ws.Range("C24").Value = "Passive"
$ws.Range("D24").Value = "Passivity of one wrestler"

# Remove the obsolete rows that used to sit right after (NearFall 2pts/4pts, Stalling,
# Technical Violation, Riding Time, Leg Attacks, Takedowns, Penalties) - the rows below
# shift up to take their place.
$ws.Range("A25:A32").EntireRow.Delete()

# Leave selection on D7, matching the saved view state.
$ws.Range("D7").Select()
